$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88; this shifts the previous rows 88-129 down to 89-130,
# automatically growing the sheet dimension to A1:R130.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new weekly record.
$ws.Cells.Item(88, 1).Value = 11
$ws.Cells.Item(88, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(88, 3).Value = "Bíobío"
$ws.Cells.Item(88, 4).Value = 44806
$ws.Cells.Item(88, 5).Value = 8
$ws.Cells.Item(88, 6).Value = 100112021
$ws.Cells.Item(88, 7).Value = "Ají"
$ws.Cells.Item(88, 8).Value = "Inferno"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 25
$ws.Cells.Item(88, 11).Value = 18000
$ws.Cells.Item(88, 12).Value = 19000
$ws.Cells.Item(88, 13).Value = 18600
$ws.Cells.Item(88, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(88, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(88, 16).Value = 1550
$ws.Cells.Item(88, 17).Value = 12
$ws.Cells.Item(88, 18).Value = "Hortaliza"
